# Applies stock-report recalculation updates described in the commit diff.
# Each block corresponds to a worksheet row whose quantity/value cells (and,
# in a few cases, swapped batch rows) were updated; Sub Total / Grand Total
# rows are updated to match the recalculated column sums.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 90
$ws.Range("F90").Value = 89
$ws.Range("G90").Value = 12010.55

# Row 92
$ws.Range("F92").Value = 109
$ws.Range("G92").Value = 12201.46

# Row 114
$ws.Range("B114").Value = 269974.58

# Row 143
$ws.Range("F143").Value = 9
$ws.Range("G143").Value = 874.8

# Row 152
$ws.Range("B152").Value = 22071.56

# Row 163
$ws.Range("B163").Value = 57552
$ws.Range("E163").Value = 136.86
$ws.Range("F163").Value = -5
$ws.Range("G163").Value = -603.45

# Row 164
$ws.Range("B164").Value = 64329
$ws.Range("E164").Value = 128.32
$ws.Range("F164").Value = 3
$ws.Range("G164").Value = 362.07

# Row 193
$ws.Range("F193").Value = 305
$ws.Range("G193").Value = 19764

# Row 195
$ws.Range("F195").Value = 126
$ws.Range("G195").Value = 10946.88

# Row 196
$ws.Range("F196").Value = 100
$ws.Range("G196").Value = 8829

# Row 200
$ws.Range("B200").Value = 48802.5

# Row 217
$ws.Range("F217").Value = 5
$ws.Range("G217").Value = 956.65

# Row 219
$ws.Range("F219").Value = 6
$ws.Range("G219").Value = 4032.24

# Row 222
$ws.Range("B222").Value = 52931.6

# Row 246
$ws.Range("B246").Value = 64973
$ws.Range("E246").Value = 35.4
$ws.Range("F246").Value = 65
$ws.Range("G246").Value = 2164.5

# Row 247
$ws.Range("B247").Value = 48706
$ws.Range("E247").Value = 39.8
$ws.Range("F247").Value = -144
$ws.Range("G247").Value = -4795.2

# Row 256
$ws.Range("F256").Value = 1
$ws.Range("G256").Value = 20.21

# Row 274
$ws.Range("B274").Value = 91859.85000000001

# Row 292
$ws.Range("B292").Value = 55373
$ws.Range("E292").Value = 163.62
$ws.Range("F292").Value = -94
$ws.Range("G292").Value = -13562.32

# Row 293
$ws.Range("B293").Value = 63520
$ws.Range("E293").Value = 153.4
$ws.Range("F293").Value = 74
$ws.Range("G293").Value = 10676.72

# Row 294
$ws.Range("B294").Value = 63571
$ws.Range("E294").Value = 152.53
$ws.Range("F294").Value = 5
$ws.Range("G294").Value = 717.4

# Row 295
$ws.Range("B295").Value = 57802
$ws.Range("E295").Value = 162.71
$ws.Range("F295").Value = -79
$ws.Range("G295").Value = -11334.92

# Row 299
$ws.Range("B299").Value = 55356
$ws.Range("E299").Value = 54.04
$ws.Range("F299").Value = -158
$ws.Range("G299").Value = -7527.12

# Row 300
$ws.Range("B300").Value = 63510
$ws.Range("E300").Value = 50.66
$ws.Range("F300").Value = 145
$ws.Range("G300").Value = 6907.8

# Row 313
$ws.Range("F313").Value = 6
$ws.Range("G313").Value = 728.16

# Row 328
$ws.Range("F328").Value = 841
$ws.Range("G328").Value = 17686.23

# Row 339
$ws.Range("B339").Value = 317170.24

# Row 378
$ws.Range("F378").Value = 9
$ws.Range("G378").Value = 667.4400000000001

# Row 395
$ws.Range("B395").Value = 250514.72

# Row 424
$ws.Range("F424").Value = 80
$ws.Range("G424").Value = 2379.2

# Row 426
$ws.Range("F426").Value = 110
$ws.Range("G426").Value = 10626

# Row 430
$ws.Range("B430").Value = 48221.91

# Row 444
$ws.Range("F444").Value = 61
$ws.Range("G444").Value = 4390.17

# Row 448
$ws.Range("B448").Value = 41245.89

# Row 455
$ws.Range("F455").Value = 37
$ws.Range("G455").Value = 5076.03

# Row 460
$ws.Range("B460").Value = 47619.8

# Row 468
$ws.Range("F468").Value = 113
$ws.Range("G468").Value = 1485.95

# Row 470
$ws.Range("F470").Value = 164
$ws.Range("G470").Value = 2100.84

# Row 477
$ws.Range("F477").Value = 120
$ws.Range("G477").Value = 2367.6

# Row 484
$ws.Range("F484").Value = 564
$ws.Range("G484").Value = 3660.36

# Row 492
$ws.Range("B492").Value = -1734.69

# Row 575
$ws.Range("F575").Value = 57
$ws.Range("G575").Value = 1511.07

# Row 578
$ws.Range("F578").Value = 60
$ws.Range("G578").Value = 6634.2

# Row 582
$ws.Range("B582").Value = 21458.82

# Row 647
$ws.Range("F647").Value = 0
$ws.Range("G647").Value = 0

# Row 650
$ws.Range("F650").Value = 354
$ws.Range("G650").Value = 28454.52

# Row 651
$ws.Range("B651").Value = 37478.92

# Row 701
$ws.Range("F701").Value = 187
$ws.Range("G701").Value = 26765.31

# Row 702
$ws.Range("F702").Value = 45
$ws.Range("G702").Value = 3670.2

# Row 705
$ws.Range("F705").Value = 70
$ws.Range("G705").Value = 5297.6

# Row 713
$ws.Range("F713").Value = 363
$ws.Range("G713").Value = 49008.63

# Row 716
$ws.Range("B716").Value = 150760.86

# Row 724
$ws.Range("F724").Value = 29
$ws.Range("G724").Value = 4361.89

# Row 741
$ws.Range("F741").Value = 124
$ws.Range("G741").Value = 29983.2

# Row 743
$ws.Range("B743").Value = 85815.34

# Row 774
$ws.Range("F774").Value = 182
$ws.Range("G774").Value = 23401.56

# Row 775
$ws.Range("B775").Value = 801251.17

# Row 778
$ws.Range("F778").Value = 111
$ws.Range("G778").Value = 16207.11

# Row 792
$ws.Range("B792").Value = 80285.32000000001

# Row 793
$ws.Range("B793").Value = 3029630.46

# Row 794
$ws.Range("B794").Value = 3029630.46
